# Update the "想去人数" (F column) values across all sheets to the
# freshly-scraped counts (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 257
$ws.Range("F3").Value = 851
$ws.Range("F5").Value = 783
$ws.Range("F7").Value = 77
$ws.Range("F8").Value = 1445
$ws.Range("F9").Value = 37823
$ws.Range("F10").Value = 7750
$ws.Range("F11").Value = 138
$ws.Range("F12").Value = 454
$ws.Range("F13").Value = 634
$ws.Range("F14").Value = 501
$ws.Range("F15").Value = 53
$ws.Range("F17").Value = 154
$ws.Range("F18").Value = 540
$ws.Range("F19").Value = 18
$ws.Range("F21").Value = 500
$ws.Range("F22").Value = 163
$ws.Range("F23").Value = 916
$ws.Range("F24").Value = 25
$ws.Range("F25").Value = 468
$ws.Range("F27").Value = 465
$ws.Range("F28").Value = 501
$ws.Range("F29").Value = 30
$ws.Range("F30").Value = 305
$ws.Range("F31").Value = 93
$ws.Range("F32").Value = 779
$ws.Range("F34").Value = 153
$ws.Range("F35").Value = 165
$ws.Range("F36").Value = 851
$ws.Range("F37").Value = 151
$ws.Range("F38").Value = 41
$ws.Range("F39").Value = 898
$ws.Range("F42").Value = 41

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 179
$ws.Range("F4").Value = 311
$ws.Range("F5").Value = 4346
$ws.Range("F7").Value = 262
$ws.Range("F10").Value = 66
$ws.Range("F11").Value = 66
$ws.Range("F12").Value = 6
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 4342
$ws.Range("F18").Value = 13

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1620
$ws.Range("F3").Value = 398

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1620
$ws.Range("F3").Value = 398
$ws.Range("F4").Value = 257
$ws.Range("F5").Value = 851
$ws.Range("F7").Value = 783
$ws.Range("F9").Value = 1445
$ws.Range("F10").Value = 37823
$ws.Range("F12").Value = 311
$ws.Range("F14").Value = 262
$ws.Range("F15").Value = 6
$ws.Range("F17").Value = 7750
$ws.Range("F19").Value = 454
$ws.Range("F20").Value = 66
$ws.Range("F21").Value = 634
$ws.Range("F22").Value = 501
$ws.Range("F23").Value = 66
$ws.Range("F24").Value = 66
$ws.Range("F25").Value = 154
$ws.Range("F26").Value = 540
$ws.Range("F27").Value = 18
$ws.Range("F28").Value = 43
$ws.Range("F29").Value = 500
$ws.Range("F30").Value = 163
$ws.Range("F31").Value = 916
$ws.Range("F32").Value = 25
$ws.Range("F33").Value = 468
$ws.Range("F35").Value = 465
$ws.Range("F36").Value = 501
$ws.Range("F37").Value = 30
$ws.Range("F38").Value = 93
$ws.Range("F39").Value = 779
$ws.Range("F41").Value = 328
$ws.Range("F42").Value = 153
$ws.Range("F43").Value = 165
$ws.Range("F44").Value = 151
$ws.Range("F45").Value = 898
$ws.Range("F47").Value = 63
$ws.Range("F48").Value = 13
$ws.Range("F49").Value = 41
